# Document sprint 11 #105
# Adds a new "Sprint 11" block (rows 155-166) mirroring the existing
# Sprint 0..10 blocks, with its own table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Copy the Sprint 10 block (rows 141:152, including formatting, merged
#    cells and the header row) down to rows 155:166 to replicate styles.
$src = $ws.Range("B141:G152")
$dst = $ws.Range("B155:G166")
$src.Copy($dst)

# 2. Fix the title text to "Sprint 11"
$ws.Range("B155").Value = "Sprint 11"

# 3. Re-enter the formulas lost during the value-only copy (E and G columns)
$ws.Range("E158").Formula = "=D158*C158"
$ws.Range("E159:E165").Formula = "=D159*C159"

$ws.Range("G158").Formula = "=F158*C158"
$ws.Range("G159:G165").Formula = "=F159*C159"

$ws.Range("E166").Formula = "=SUM(E158:E165)/60"
$ws.Range("G166").Formula = "=SUM(G158:G165)/60"

# 4. Update the "Min." (D) and "Real" (F) input columns with the Sprint 11 data
$ws.Range("D158").Value = 1
$ws.Range("D159").Value = 2
$ws.Range("D160").Value = 1
$ws.Range("D161").Value = 2
$ws.Range("D162").Value = 1
$ws.Range("D163").Value = 0
$ws.Range("D164").Value = 0
$ws.Range("D165").Value = 0

$ws.Range("F158").Value = 1
$ws.Range("F159").Value = 1
$ws.Range("F160").Value = 2
$ws.Range("F161").Value = 1
$ws.Range("F162").Value = 2
$ws.Range("F163").Value = 0
$ws.Range("F164").Value = 0
$ws.Range("F165").Value = 0

# 5. Turn the new range into a table, matching the other sprint tables
$rng = $ws.Range("B157:G166")
$lo = $ws.ListObjects.Add(1, $rng, $null, 1)
$lo.Name = "Tabla1810112345671213"
$lo.TableStyle = "TableStyleMedium16"
$lo.ShowTableStyleRowStripes = $false
$lo.ShowTableStyleColumnStripes = $true
$lo.TotalsRowShown = $false

# 6. Update views so the new block is visible / selected, matching the author
$ws.Range("E166").Select()
